# Update cryptos list: refresh Price / Volume(1h) values, and for row 12/13
# swap Chainlink and WrappedliquidstakedEther2.0 (with their updated figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'37.706.84"
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = "'2.030.03"
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').Value = "'227.19"
$ws.Range('E5').Value = '  -1.67%  '
$ws.Range('D6').Value = "'0.608"
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('D7').Value = "'59.85"
$ws.Range('E7').Value = '  -2.40%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -3.41%  '
$ws.Range('D10').Value = "'0.0827"
$ws.Range('E10').Value = '  +2.21%  '
$ws.Range('D11').Value = "'0.104"
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('B12').Value = 'Chainlink'
$ws.Range('C12').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D12').Value = "'14.56"
$ws.Range('E12').Value = '  -2.94%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = "'2.330.23"
$ws.Range('E13').Value = '  -1.87%  '
$ws.Range('D14').Value = "'21.02"
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = "'5.19"
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('D17').Value = "'2.029.24"
$ws.Range('E17').Value = '  -2.45%  '
$ws.Range('D18').Value = "'37.681.29"
$ws.Range('E18').Value = '  -1.40%  '
$ws.Range('D19').Value = "'69.47"
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').Value = "'5.90"
$ws.Range('E20').Value = '  -6.58%  '
$ws.Range('D22').Value = "'223.75"
$ws.Range('E22').Value = '  -0.93%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  -2.76%  '
$ws.Range('E25').Value = '  -0.96%  '
$ws.Range('D26').Value = "'168.07"
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D27').Value = "'9.36"
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('E28').Value = '  -3.55%  '
$ws.Range('D29').Value = "'18.76"
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('D30').Value = "'1.26"
$ws.Range('E30').Value = '  -5.55%  '
$ws.Range('E31').Value = '  +0.47%  '
$ws.Range('E32').Value = '  +7.68%  '
$ws.Range('E33').Value = '  -4.54%  '
$ws.Range('D34').Value = "'0.0604"
$ws.Range('E34').Value = '  -1.11%  '
$ws.Range('E35').Value = '  -4.66%  '
$ws.Range('D36').Value = "'6.45"
$ws.Range('E36').Value = '  +2.68%  '
$ws.Range('E37').Value = '  +0.20%  '
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = "'18.04"
$ws.Range('E40').Value = '  +4.09%  '
$ws.Range('D41').Value = "'1.538.40"
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('E42').Value = '  -1.52%  '
$ws.Range('D43').Value = "'95.51"
$ws.Range('E43').Value = '  -3.36%  '
$ws.Range('D44').Value = "'2.79"
$ws.Range('E44').Value = '  -2.92%  '
$ws.Range('D45').Value = "'0.0908"
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').Value = "'4.10"
$ws.Range('E46').Value = '  +1.33%  '
$ws.Range('E47').Value = '  -3.23%  '
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('D51').Value = "'2.218.99"
$ws.Range('E51').Value = '  -1.88%  '
